# Applies the "Updated symbol list" commit: refreshes prices and, for a
# handful of rows, rotates the Coin/Link/Volume columns to reflect new
# rankings, exactly as captured in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells hold numeric-looking text (e.g. "248.63") that must stay a
# text value (matching the workbook's inlineStr cells), not be coerced to
# a number by Excel's normal input parsing. Forcing the cell to Text
# format before the write keeps it literal; ClearFormats() afterwards
# drops the now-unneeded number format so no stray style sticks to the
# cell.
function Set-TextCell($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

function Set-Cell($addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- Simple price-only updates -------------------------------------------
Set-TextCell "D2"  "248.63"
Set-TextCell "D3"  "24.07"
Set-TextCell "D4"  "5.820"
Set-TextCell "D6"  "3.434"
Set-TextCell "D7"  "6.513"
Set-TextCell "D8"  "1.331"
Set-TextCell "D9"  "0.7965"

# --- Rows 10-18: coin ranking rotated up by one row -----------------------
Set-Cell "B10" "WazirX"
Set-Cell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1473"
Set-Cell "E10" "9WazirXWRX"

Set-Cell "B11" "MandalaExchangeToken"
Set-Cell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.07737"
Set-Cell "E11" "10MandalaExchangeTokenMDX"

Set-Cell "B12" "LiechtensteinCryptoassetsExchange"
Set-Cell "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D12" "0.03293"
Set-Cell "E12" "11LiechtensteinCryptoassetsExchangeLCX"

Set-Cell "B13" "BitrueCoin"
Set-Cell "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.03015"
Set-Cell "E13" "12BitrueCoinBTR"

Set-Cell "B14" "BitMartToken"
Set-Cell "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D14" "0.09233"
Set-Cell "E14" "13BitMartTokenBMX"

Set-Cell "B15" "MCDex"
Set-Cell "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D15" "3.570"
Set-Cell "E15" "14MCDexMCB"

Set-Cell "B16" "BitForexToken"
Set-Cell "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D16" "0.001660"
Set-Cell "E16" "15BitForexTokenBF"

Set-Cell "B17" "CoinExToken"
Set-Cell "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D17" "0.04763"
Set-Cell "E17" "16CoinExTokenCET"

Set-Cell "B18" "One"
Set-Cell "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D18" "0.0006038"
Set-Cell "E18" "17OneONE"

# --- More simple price-only updates ---------------------------------------
Set-TextCell "D19" "0.006226"
Set-TextCell "D20" "0.005528"
Set-TextCell "D21" "0.001068"
Set-TextCell "D22" "0.0001501"
Set-TextCell "D23" "3.694"
Set-TextCell "D25" "0.3351"
Set-TextCell "D27" "0.0006253"
Set-TextCell "D40" "0.04375"
Set-TextCell "D41" "0.007019"

# --- Rows 42-43: CEJI / BKEXToken swap ranking -----------------------------
Set-Cell "B42" "BKEXToken"
Set-Cell "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D42" "0.1061"
Set-Cell "E42" "41BKEXTokenBKK"

Set-Cell "B43" "CEJI"
Set-Cell "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D43" "0.003214"
Set-Cell "E43" "42CEJICEJI"

# --- Remaining updates ------------------------------------------------------
Set-TextCell "D44" "0.009659"
Set-TextCell "D45" "0.002463"
Set-Cell "E45" "44ACDXExchangeACXTBestin24h"
Set-TextCell "D46" "0.00005895"
Set-TextCell "D48" "0.9916"
Set-TextCell "D49" "0.1110"
Set-Cell "E49" "48BOLOBOLOWorstin24h"
Set-TextCell "D50" "0.00002103"
Set-Cell "E51" "50SpecialPowerGoldSPG"
